$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Determine the last used row on the sheet (data runs from row 2 to this row).
$lastRow = $ws.UsedRange.Rows.Count - 1

# Update the "Förändrad" (changed) date column C for every data row (2..lastRow)
# from 45172 to 45175 (re-scrape timestamp bump).
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45172) {
        $cell.Value = 45175
    }
}

# Row 3 ("A 64788-2019") gained a new species finding: "Vaddporing".
# That bumps NT (J), Rödlistade (O) and Alla arter (Q) counts by one each,
# and the species name is inserted into the species list (R) right after
# "Ullticka" and before "Vedskivlav".
$ws.Cells.Item(3, 10).Value = 10   # J3 NT
$ws.Cells.Item(3, 15).Value = 10   # O3 Rödlistade
$ws.Cells.Item(3, 17).Value = 13   # Q3 Alla arter

$r3 = $ws.Cells.Item(3, 18)
$oldText = $r3.Value2
$newText = $oldText -replace "Ullticka`r`n", "Ullticka`r`nVaddporing`r`n"
$r3.Value = $newText
